# Apply edits described in the commit diff:
# 1) Update the date/weekday header line.
# 2) Update the 25 division problems in the practice table.

$d = $word.ActiveDocument

# --- 1) Update header date/weekday -----------------------------------
$d.Content.Find.Execute("2025-01-15 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-01-16 Thursday", 2) | Out-Null

# --- 2) Update the division problems in the table ---------------------
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "37÷4="
$t.Cell(1,2).Range.Text  = "23÷4="
$t.Cell(1,3).Range.Text  = "87÷7="
$t.Cell(1,4).Range.Text  = "30÷8="
$t.Cell(1,5).Range.Text  = "69÷8="

$t.Cell(5,1).Range.Text  = "46÷5="
$t.Cell(5,2).Range.Text  = "79÷3="
$t.Cell(5,3).Range.Text  = "50÷9="
$t.Cell(5,4).Range.Text  = "11÷5="
$t.Cell(5,5).Range.Text  = "76÷3="

$t.Cell(9,1).Range.Text  = "24÷3="
$t.Cell(9,2).Range.Text  = "42÷9="
$t.Cell(9,3).Range.Text  = "53÷2="
$t.Cell(9,4).Range.Text  = "63÷3="
$t.Cell(9,5).Range.Text  = "90÷3="

$t.Cell(13,1).Range.Text = "48÷9="
$t.Cell(13,2).Range.Text = "36÷5="
$t.Cell(13,3).Range.Text = "91÷8="
$t.Cell(13,4).Range.Text = "12÷2="
$t.Cell(13,5).Range.Text = "61÷8="

$t.Cell(17,1).Range.Text = "29÷9="
$t.Cell(17,2).Range.Text = "75÷6="
$t.Cell(17,3).Range.Text = "10÷9="
$t.Cell(17,4).Range.Text = "32÷5="
$t.Cell(17,5).Range.Text = "72÷6="

Write-Host "Edits applied."
